# Add a new "Customer defect name" column before the existing "Set Model"
# column (old column K), shifting K:R -> L:S, and update the dependent
# data-validation references / selected cell accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at K; Excel automatically shifts the existing K:R
# columns (and their cell contents/number formats) one slot to the right.
$ws.Columns("K:K").Insert()

# New column K should have the same width as column J (the XML stores both
# as width 29). Re-reading ColumnWidth and writing it back round-trips to
# the same underlying "29" the source column already serializes as.
$ws.Columns("K:K").ColumnWidth = $ws.Columns("J:J").ColumnWidth

# Give the new header cell K1 the same visual style as the other header
# cells (fill + border), set its text, then recolor the font red to flag it
# as the new "customer defect name" column.
$ws.Range("J1").Copy()
$ws.Range("K1").PasteSpecial(-4122)
$ws.Range("K1").Value = "Customer defect name"
$ws.Range("K1").Font.Color = 255

# K2 stays empty but carries the same bordered/number-format style as the
# other date cell in row 2 (matches style id used by J2).
$ws.Range("J2").Copy()
$ws.Range("K2").PasteSpecial(-4122)
$ws.Range("K2").ClearContents()

# The data-validation list sources (previously columns P:R, now Q:S) need
# their referenced ranges updated to follow the column shift. Editing
# Formula1 in place (rather than delete + re-add) preserves the existing
# sqref grouping, e.g. "M2 M5" stays combined in one <dataValidation>.
$ws.Range("D2").Validation.Formula1 = "=`$Q`$1:`$Q`$2"
$ws.Range("E2").Validation.Formula1 = "=`$R`$1:`$R`$2"
$ws.Range("M2").Validation.Formula1 = "=`$S`$1:`$S`$2"

# Move the active selection to match the saved view state.
$ws.Range("I7").Select()
